$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Columns("G").Insert()

$rng2 = $ws.Range("G6")
$rng2.Style = "Hyperlink"
$rng2.Font.Name = "Arial"
$rng2.Font.Size = 10
$rng2.Font.Underline = $false
$rng2.Font.ThemeColor = 1
Write-Host "done"
